$wb = $excel.ActiveWorkbook

# Remove the extra "Sheet1" worksheet (and its embedded Table1) entirely
$excel.DisplayAlerts = $false
$null = $wb.Worksheets.Item("Sheet1").Delete()
$excel.DisplayAlerts = $true

$ws = $wb.Worksheets.Item("master-template_type")

# Append the new master-data rows (92-121) to the main sheet
$ws.Cells.Item(92,1).Value = "RPR_UIN_CARD_TEMPLATE"
$ws.Cells.Item(92,2).Value = "UIN card template"
$ws.Cells.Item(92,3).Value = "eng"
$ws.Cells.Item(92,4).Value = $true
$ws.Cells.Item(92,5).Value = "superadmin"
$ws.Cells.Item(92,6).Value = "now()"

$ws.Cells.Item(93,1).Value = "RPR_UIN_CARD_TEMPLATE"
$ws.Cells.Item(93,2).Value = "قالب بطاقة UIN"
$ws.Cells.Item(93,3).Value = "ara"
$ws.Cells.Item(93,4).Value = $true
$ws.Cells.Item(93,5).Value = "superadmin"
$ws.Cells.Item(93,6).Value = "now()"

$ws.Cells.Item(94,1).Value = "RPR_UIN_CARD_TEMPLATE"
$ws.Cells.Item(94,2).Value = "Modèle de carte UIN"
$ws.Cells.Item(94,3).Value = "fra"
$ws.Cells.Item(94,4).Value = $true
$ws.Cells.Item(94,5).Value = "superadmin"
$ws.Cells.Item(94,6).Value = "now()"

$ws.Cells.Item(95,1).Value = "RPR_UIN_DEAC_SMS"
$ws.Cells.Item(95,2).Value = "Template for UIN Deactivation SMS"
$ws.Cells.Item(95,3).Value = "eng"
$ws.Cells.Item(95,4).Value = $true
$ws.Cells.Item(95,5).Value = "superadmin"
$ws.Cells.Item(95,6).Value = "now()"

$ws.Cells.Item(96,1).Value = "RPR_UIN_DEAC_SMS"
$ws.Cells.Item(96,2).Value = "قالب لتعطيل UIN SMS"
$ws.Cells.Item(96,3).Value = "ara"
$ws.Cells.Item(96,4).Value = $true
$ws.Cells.Item(96,5).Value = "superadmin"
$ws.Cells.Item(96,6).Value = "now()"

$ws.Cells.Item(97,1).Value = "RPR_UIN_DEAC_SMS"
$ws.Cells.Item(97,2).Value = "Modèle pour SMS de désactivation UIN"
$ws.Cells.Item(97,3).Value = "fra"
$ws.Cells.Item(97,4).Value = $true
$ws.Cells.Item(97,5).Value = "superadmin"
$ws.Cells.Item(97,6).Value = "now()"

$ws.Cells.Item(98,1).Value = "RPR_UIN_DEAC_EMAIL"
$ws.Cells.Item(98,2).Value = "Template for UIN Deactivation Email"
$ws.Cells.Item(98,3).Value = "eng"
$ws.Cells.Item(98,4).Value = $true
$ws.Cells.Item(98,5).Value = "superadmin"
$ws.Cells.Item(98,6).Value = "now()"

$ws.Cells.Item(99,1).Value = "RPR_UIN_DEAC_EMAIL"
$ws.Cells.Item(99,2).Value = "قالب لإلغاء تنشيط البريد"
$ws.Cells.Item(99,3).Value = "ara"
$ws.Cells.Item(99,4).Value = $true
$ws.Cells.Item(99,5).Value = "superadmin"
$ws.Cells.Item(99,6).Value = "now()"

$ws.Cells.Item(100,1).Value = "RPR_UIN_DEAC_EMAIL"
$ws.Cells.Item(100,2).Value = "Modèle pour Email de désactivation UIN"
$ws.Cells.Item(100,3).Value = "fra"
$ws.Cells.Item(100,4).Value = $true
$ws.Cells.Item(100,5).Value = "superadmin"
$ws.Cells.Item(100,6).Value = "now()"

$ws.Cells.Item(101,1).Value = "RPR_UIN_REAC_SMS"
$ws.Cells.Item(101,2).Value = "Template for UIN Reactivate SMS"
$ws.Cells.Item(101,3).Value = "eng"
$ws.Cells.Item(101,4).Value = $true
$ws.Cells.Item(101,5).Value = "superadmin"
$ws.Cells.Item(101,6).Value = "now()"

$ws.Cells.Item(102,1).Value = "RPR_UIN_REAC_SMS"
$ws.Cells.Item(102,2).Value = "قالب لـ UIN تنشيط SMS"
$ws.Cells.Item(102,3).Value = "ara"
$ws.Cells.Item(102,4).Value = $true
$ws.Cells.Item(102,5).Value = "superadmin"
$ws.Cells.Item(102,6).Value = "now()"

$ws.Cells.Item(103,1).Value = "RPR_UIN_REAC_SMS"
$ws.Cells.Item(103,2).Value = "Modèle pour UIN Réactiver SMS"
$ws.Cells.Item(103,3).Value = "fra"
$ws.Cells.Item(103,4).Value = $true
$ws.Cells.Item(103,5).Value = "superadmin"
$ws.Cells.Item(103,6).Value = "now()"

$ws.Cells.Item(104,1).Value = "RPR_UIN_REAC_EMAIL"
$ws.Cells.Item(104,2).Value = "Template for UIN Reactivate Email"
$ws.Cells.Item(104,3).Value = "eng"
$ws.Cells.Item(104,4).Value = $true
$ws.Cells.Item(104,5).Value = "superadmin"
$ws.Cells.Item(104,6).Value = "now()"

$ws.Cells.Item(105,1).Value = "RPR_UIN_REAC_EMAIL"
$ws.Cells.Item(105,2).Value = "قالب لـ UIN تنشيط البريد"
$ws.Cells.Item(105,3).Value = "ara"
$ws.Cells.Item(105,4).Value = $true
$ws.Cells.Item(105,5).Value = "superadmin"
$ws.Cells.Item(105,6).Value = "now()"

$ws.Cells.Item(106,1).Value = "RPR_UIN_REAC_EMAIL"
$ws.Cells.Item(106,2).Value = "Modèle pour UIN Réactiver Email"
$ws.Cells.Item(106,3).Value = "fra"
$ws.Cells.Item(106,4).Value = $true
$ws.Cells.Item(106,5).Value = "superadmin"
$ws.Cells.Item(106,6).Value = "now()"

$ws.Cells.Item(107,1).Value = "reg-sms-notification"
$ws.Cells.Item(107,2).Value = "Registration Acknowledgement Template"
$ws.Cells.Item(107,3).Value = "eng"
$ws.Cells.Item(107,4).Value = $true
$ws.Cells.Item(107,5).Value = "superadmin"
$ws.Cells.Item(107,6).Value = "now()"

$ws.Cells.Item(108,1).Value = "reg-sms-notification"
$ws.Cells.Item(108,2).Value = "نموذج شكر التسجيل"
$ws.Cells.Item(108,3).Value = "ara"
$ws.Cells.Item(108,4).Value = $true
$ws.Cells.Item(108,5).Value = "superadmin"
$ws.Cells.Item(108,6).Value = "now()"

$ws.Cells.Item(109,1).Value = "reg-sms-notification"
$ws.Cells.Item(109,2).Value = "accusé de réception"
$ws.Cells.Item(109,3).Value = "fra"
$ws.Cells.Item(109,4).Value = $true
$ws.Cells.Item(109,5).Value = "superadmin"
$ws.Cells.Item(109,6).Value = "now()"

$ws.Cells.Item(110,1).Value = "reg-email-notification"
$ws.Cells.Item(110,2).Value = "Registration Acknowledgement Template"
$ws.Cells.Item(110,3).Value = "eng"
$ws.Cells.Item(110,4).Value = $true
$ws.Cells.Item(110,5).Value = "superadmin"
$ws.Cells.Item(110,6).Value = "now()"

$ws.Cells.Item(111,1).Value = "reg-email-notification"
$ws.Cells.Item(111,2).Value = "نموذج شكر التسجيل"
$ws.Cells.Item(111,3).Value = "ara"
$ws.Cells.Item(111,4).Value = $true
$ws.Cells.Item(111,5).Value = "superadmin"
$ws.Cells.Item(111,6).Value = "now()"

$ws.Cells.Item(112,1).Value = "reg-email-notification"
$ws.Cells.Item(112,2).Value = "accusé de réception"
$ws.Cells.Item(112,3).Value = "fra"
$ws.Cells.Item(112,4).Value = $true
$ws.Cells.Item(112,5).Value = "superadmin"
$ws.Cells.Item(112,6).Value = "now()"

$ws.Cells.Item(113,1).Value = "reg-ack-template-part1"
$ws.Cells.Item(113,2).Value = "Registration Acknowledgement Template - Part 1"
$ws.Cells.Item(113,3).Value = "eng"
$ws.Cells.Item(113,4).Value = $true
$ws.Cells.Item(113,5).Value = "superadmin"
$ws.Cells.Item(113,6).Value = "now()"

$ws.Cells.Item(114,1).Value = "reg-ack-template-part2"
$ws.Cells.Item(114,2).Value = "نموذج شكر التسجيل"
$ws.Cells.Item(114,3).Value = "ara"
$ws.Cells.Item(114,4).Value = $true
$ws.Cells.Item(114,5).Value = "superadmin"
$ws.Cells.Item(114,6).Value = "now()"

$ws.Cells.Item(115,1).Value = "reg-ack-template-part3"
$ws.Cells.Item(115,2).Value = "accusé de réception"
$ws.Cells.Item(115,3).Value = "fra"
$ws.Cells.Item(115,4).Value = $true
$ws.Cells.Item(115,5).Value = "superadmin"
$ws.Cells.Item(115,6).Value = "now()"

$ws.Cells.Item(116,1).Value = "reg-ack-template-part2"
$ws.Cells.Item(116,2).Value = "Registration Acknowledgement Template - Part 2"
$ws.Cells.Item(116,3).Value = "eng"
$ws.Cells.Item(116,4).Value = $true
$ws.Cells.Item(116,5).Value = "superadmin"
$ws.Cells.Item(116,6).Value = "now()"

$ws.Cells.Item(117,1).Value = "reg-ack-template-part3"
$ws.Cells.Item(117,2).Value = "نموذج شكر التسجيل"
$ws.Cells.Item(117,3).Value = "ara"
$ws.Cells.Item(117,4).Value = $true
$ws.Cells.Item(117,5).Value = "superadmin"
$ws.Cells.Item(117,6).Value = "now()"

$ws.Cells.Item(118,1).Value = "reg-ack-template-part4"
$ws.Cells.Item(118,2).Value = "accusé de réception"
$ws.Cells.Item(118,3).Value = "fra"
$ws.Cells.Item(118,4).Value = $true
$ws.Cells.Item(118,5).Value = "superadmin"
$ws.Cells.Item(118,6).Value = "now()"

$ws.Cells.Item(119,1).Value = "reg-ack-template-part3"
$ws.Cells.Item(119,2).Value = "Registration Acknowledgement Template - Part 3"
$ws.Cells.Item(119,3).Value = "eng"
$ws.Cells.Item(119,4).Value = $true
$ws.Cells.Item(119,5).Value = "superadmin"
$ws.Cells.Item(119,6).Value = "now()"

$ws.Cells.Item(120,1).Value = "reg-ack-template-part4"
$ws.Cells.Item(120,2).Value = "نموذج شكر التسجيل"
$ws.Cells.Item(120,3).Value = "ara"
$ws.Cells.Item(120,4).Value = $true
$ws.Cells.Item(120,5).Value = "superadmin"
$ws.Cells.Item(120,6).Value = "now()"

$ws.Cells.Item(121,1).Value = "reg-ack-template-part5"
$ws.Cells.Item(121,2).Value = "accusé de réception"
$ws.Cells.Item(121,3).Value = "fra"
$ws.Cells.Item(121,4).Value = $true
$ws.Cells.Item(121,5).Value = "superadmin"
$ws.Cells.Item(121,6).Value = "now()"

# Move the active selection below the newly added data, matching the saved workbook state
$null = $ws.Range("A122:XFD1048576").Select()

